$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for the two new columns (P=16, Q=17) in row 1,
# matching the bold/centered/bordered style already used by the rest
# of the header row (e.g. O1).
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108
$ws.Range("P1:Q1").VerticalAlignment = -4160
$ws.Range("P1:Q1").Borders.LineStyle = 1

# For every data row, swap the I/K and M/O column values and populate the
# two new columns (P, Q) with the value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column
}
